$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the existing used range (A1:AD19) completely (values + formatting) ---
$ws.Range("A1:AD19").Clear()

# --- Build the full new table (23 rows x 20 cols, A1:T23) ---
$data = New-Object 'object[,]' 23,20
$data[0,1] = 0
$data[0,2] = 1
$data[0,3] = 2
$data[0,4] = 3
$data[0,5] = 4
$data[0,6] = 5
$data[0,7] = 6
$data[0,8] = 7
$data[0,9] = 8
$data[0,10] = 9
$data[0,11] = 10
$data[0,12] = 11
$data[0,13] = 12
$data[0,14] = 13
$data[0,15] = 14
$data[0,16] = 15
$data[0,17] = 16
$data[0,18] = 17
$data[0,19] = 18
$data[1,0] = 0
$data[1,1] = "HKL"
$data[1,2] = "[2, 0, 0]"
$data[1,3] = "[2, 2, 0]"
$data[1,4] = "[4, 0, 0]"
$data[1,5] = "[2, 1, 1]"
$data[1,6] = "[3, 2, 1]"
$data[1,7] = "[3, 1, 0]"
$data[1,8] = "[2, 2, 2]"
$data[1,9] = "[1, 1, 0]"
$data[1,10] = "1Pair-A"
$data[1,11] = "1Pair-B"
$data[1,12] = "2Pairs-A"
$data[1,13] = "2Pairs-B"
$data[1,14] = "3Pairs-A"
$data[1,15] = "3Pairs-B"
$data[1,16] = "3Pairs-C"
$data[1,17] = "4Pairs"
$data[1,18] = "5A4F"
$data[1,19] = "MaxUnique"
$data[2,0] = 1
$data[2,1] = "BT8Hex_2.5"
$data[2,2] = 1.000261405711216
$data[2,3] = 0.9999346459020385
$data[2,4] = 1.000261405711216
$data[2,5] = 0.9999346459020385
$data[2,6] = 0.9999346459020385
$data[2,7] = 1.000143771385049
$data[2,8] = 0.9998257291731821
$data[2,9] = 0.9999346459020385
$data[2,10] = 0.9999346459020385
$data[2,11] = 0.9999346459020385
$data[2,12] = 1.000098025806627
$data[2,13] = 1.000098025806627
$data[2,14] = 1.000113274332768
$data[2,15] = 1.000043565838431
$data[2,16] = 1.000043565838431
$data[2,17] = 1.000016335854333
$data[2,18] = 1.000016335854333
$data[2,19] = 1.00000580732926
$data[3,0] = 2
$data[3,1] = "BT8Hex_5"
$data[3,2] = 1.000504638976168
$data[3,3] = 0.9998738381124863
$data[3,4] = 1.000504638976168
$data[3,5] = 0.9998738381124863
$data[3,6] = 0.9998738381124863
$data[3,7] = 1.000277550320195
$data[3,8] = 0.9996635714710189
$data[3,9] = 0.9998738381124863
$data[3,10] = 0.9998738381124863
$data[3,11] = 0.9998738381124863
$data[3,12] = 1.000189238544327
$data[3,13] = 1.000189238544327
$data[3,14] = 1.00021867580295
$data[3,15] = 1.000084105067047
$data[3,16] = 1.000084105067047
$data[3,17] = 1.000031538328407
$data[3,18] = 1.000031538328407
$data[3,19] = 1.000011212517473
$data[4,0] = 3
$data[4,1] = "BT8Hex_10"
$data[4,2] = 1.00096708360052
$data[4,3] = 0.9997582274323231
$data[4,4] = 1.00096708360052
$data[4,5] = 0.9997582274323231
$data[4,6] = 0.9997582274323231
$data[4,7] = 1.000531893966553
$data[4,8] = 0.9993552796067325
$data[4,9] = 0.9997582274323231
$data[4,10] = 0.9997582274323231
$data[4,11] = 0.9997582274323231
$data[4,12] = 1.000362655516422
$data[4,13] = 1.000362655516422
$data[4,14] = 1.000419068333132
$data[4,15] = 1.000161179488389
$data[4,16] = 1.000161179488389
$data[4,17] = 1.000060441474372
$data[4,18] = 1.000060441474372
$data[4,19] = 1.000021489911796
$data[5,0] = 4
$data[5,1] = "BT8Hex_15"
$data[5,2] = 1.001418396481083
$data[5,3] = 0.9996453997972402
$data[5,4] = 1.001418396481083
$data[5,5] = 0.9996453997972402
$data[5,6] = 0.9996453997972402
$data[5,7] = 1.000780107630316
$data[5,8] = 0.9990544008062195
$data[5,9] = 0.9996453997972402
$data[5,10] = 0.9996453997972402
$data[5,11] = 0.9996453997972402
$data[5,12] = 1.000531898139162
$data[5,13] = 1.000531898139162
$data[5,14] = 1.000614634636213
$data[5,15] = 1.000236398691855
$data[5,16] = 1.000236398691855
$data[5,17] = 1.000088648968201
$data[5,18] = 1.000088648968201
$data[5,19] = 1.00003151738489
$data[6,0] = 5
$data[6,1] = "Spiral2.5"
$data[6,2] = 1.000017132146924
$data[6,3] = 0.9999957145475029
$data[6,4] = 1.000017132146924
$data[6,5] = 0.9999957145475029
$data[6,6] = 0.9999957145475029
$data[6,7] = 1.000009421414056
$data[6,8] = 0.9999885779639821
$data[6,9] = 0.9999957145475029
$data[6,10] = 0.9999957145475029
$data[6,11] = 0.9999957145475029
$data[6,12] = 1.000006423347214
$data[6,13] = 1.000006423347214
$data[6,14] = 1.000007422702828
$data[6,15] = 1.00000285374731
$data[6,16] = 1.00000285374731
$data[6,17] = 1.000001068947358
$data[6,18] = 1.000001068947358
$data[6,19] = 1.000000379194579
$data[7,0] = 6
$data[7,1] = "Spiral5"
$data[7,2] = 1.000039883468302
$data[7,3] = 0.9999900267190296
$data[7,4] = 1.000039883468302
$data[7,5] = 0.9999900267190296
$data[7,6] = 0.9999900267190296
$data[7,7] = 1.000021934678361
$data[7,8] = 0.9999734106279775
$data[7,9] = 0.9999900267190296
$data[7,10] = 0.9999900267190296
$data[7,11] = 0.9999900267190296
$data[7,12] = 1.000014955093666
$data[7,13] = 1.000014955093666
$data[7,14] = 1.000017281621898
$data[7,15] = 1.000006645635454
$data[7,16] = 1.000006645635454
$data[7,17] = 1.000002490906348
$data[7,18] = 1.000002490906348
$data[7,19] = 1.000000884821955
$data[8,0] = 7
$data[8,1] = "Spiral7.5"
$data[8,2] = 1.00005312462476
$data[8,3] = 0.9999867165035378
$data[8,4] = 1.00005312462476
$data[8,5] = 0.9999867165035378
$data[8,6] = 0.9999867165035378
$data[8,7] = 1.000029217134552
$data[8,8] = 0.9999645831053817
$data[8,9] = 0.9999867165035378
$data[8,10] = 0.9999867165035378
$data[8,11] = 0.9999867165035378
$data[8,12] = 1.000019920564149
$data[8,13] = 1.000019920564149
$data[8,14] = 1.00002301942095
$data[8,15] = 1.000008852543945
$data[8,16] = 1.000008852543945
$data[8,17] = 1.000003318533843
$data[8,18] = 1.000003318533843
$data[8,19] = 1.000001179062551
$data[9,0] = 8
$data[9,1] = "Spiral10"
$data[9,2] = 1.00011720907637
$data[9,3] = 0.9999706952003459
$data[9,4] = 1.00011720907637
$data[9,5] = 0.9999706952003459
$data[9,6] = 0.9999706952003459
$data[9,7] = 1.000064462919126
$data[9,8] = 0.9999218608965403
$data[9,9] = 0.9999706952003459
$data[9,10] = 0.9999706952003459
$data[9,11] = 0.9999706952003459
$data[9,12] = 1.000043952138358
$data[9,13] = 1.000043952138358
$data[9,14] = 1.000050789065281
$data[9,15] = 1.00001953315902
$data[9,16] = 1.00001953315902
$data[9,17] = 1.000007323669352
$data[9,18] = 1.000007323669352
$data[9,19] = 1.000002603082179
$data[10,0] = 9
$data[10,1] = "Spiral15"
$data[10,2] = 1.000186599095983
$data[10,3] = 0.9999533474123841
$data[10,4] = 1.000186599095983
$data[10,5] = 0.9999533474123841
$data[10,6] = 0.9999533474123841
$data[10,7] = 1.000102627263866
$data[10,8] = 0.9998755986670198
$data[10,9] = 0.9999533474123841
$data[10,10] = 0.9999533474123841
$data[10,11] = 0.9999533474123841
$data[10,12] = 1.000069973254183
$data[10,13] = 1.000069973254183
$data[10,14] = 1.000080857924078
$data[10,15] = 1.000031097973584
$data[10,16] = 1.000031097973584
$data[10,17] = 1.000011660333284
$data[10,18] = 1.000011660333284
$data[10,19] = 1.000004144544003
$data[11,0] = 10
$data[11,1] = "OffsetF45"
$data[11,2] = 0.9987442807092567
$data[11,3] = 1.000313923165885
$data[11,4] = 0.9987442807092567
$data[11,5] = 1.000313923165885
$data[11,6] = 1.000313923165885
$data[11,7] = 0.9993093563590919
$data[11,8] = 1.000837139764808
$data[11,9] = 1.000313923165885
$data[11,10] = 1.000313923165885
$data[11,11] = 1.000313923165885
$data[11,12] = 0.9995291019375707
$data[11,13] = 0.9995291019375707
$data[11,14] = 0.9994558534114111
$data[11,15] = 0.9997907090136753
$data[11,16] = 0.9997907090136753
$data[11,17] = 0.9999215125517276
$data[11,18] = 0.9999215125517276
$data[11,19] = 0.9999720910551351
$data[12,0] = 11
$data[12,1] = "OffsetA45"
$data[12,2] = 0.9996777471256533
$data[12,3] = 1.000080556694701
$data[12,4] = 0.9996777471256533
$data[12,5] = 1.000080556694701
$data[12,6] = 1.000080556694701
$data[12,7] = 0.9998227609416107
$data[12,8] = 1.000214833475117
$data[12,9] = 1.000080556694701
$data[12,10] = 1.000080556694701
$data[12,11] = 1.000080556694701
$data[12,12] = 0.9998791519101771
$data[12,13] = 0.9998791519101771
$data[12,14] = 0.999860354920655
$data[12,15] = 0.9999462868383517
$data[12,16] = 0.9999462868383517
$data[12,17] = 0.999979854302439
$data[12,18] = 0.999979854302439
$data[12,19] = 0.9999928352710805
$data[13,0] = 12
$data[13,1] = "OffsetFTD"
$data[13,2] = 0.9970390191027829
$data[13,3] = 1.000740231517839
$data[13,4] = 0.9970390191027829
$data[13,5] = 1.000740231517839
$data[13,6] = 1.000740231517839
$data[13,7] = 0.9983714688877133
$data[13,8] = 1.001973980505824
$data[13,9] = 1.000740231517839
$data[13,10] = 1.000740231517839
$data[13,11] = 1.000740231517839
$data[13,12] = 0.9988896253103111
$data[13,13] = 0.9988896253103111
$data[13,14] = 0.9987169065027786
$data[13,15] = 0.9995064940461539
$data[13,16] = 0.9995064940461539
$data[13,17] = 0.9998149284140752
$data[13,18] = 0.9998149284140752
$data[13,19] = 0.9999341938416396
$data[14,0] = 13
$data[14,1] = "OffsetATD"
$data[14,2] = 0.9992335485138972
$data[14,3] = 1.000191609808738
$data[14,4] = 0.9992335485138972
$data[14,5] = 1.000191609808738
$data[14,6] = 1.000191609808738
$data[14,7] = 0.9995784494664681
$data[14,8] = 1.000510962408722
$data[14,9] = 1.000191609808738
$data[14,10] = 1.000191609808738
$data[14,11] = 1.000191609808738
$data[14,12] = 0.9997125791613178
$data[14,13] = 0.9997125791613178
$data[14,14] = 0.9996678692630345
$data[14,15] = 0.9998722560437914
$data[14,16] = 0.9998722560437914
$data[14,17] = 0.9999520944850282
$data[14,18] = 0.9999520944850282
$data[14,19] = 0.9999829649692171
$data[15,0] = 14
$data[15,1] = "Holden2.5"
$data[15,2] = 1.005097787212923
$data[15,3] = 0.9987255496599942
$data[15,4] = 1.005097787212923
$data[15,5] = 0.9987255496599942
$data[15,6] = 0.9987255496599942
$data[15,7] = 1.00280377576233
$data[15,8] = 0.9966014764646669
$data[15,9] = 0.9987255496599942
$data[15,10] = 0.9987255496599942
$data[15,11] = 0.9987255496599942
$data[15,12] = 1.001911668436459
$data[15,13] = 1.001911668436459
$data[15,14] = 1.002209037545083
$data[15,15] = 1.000849628844304
$data[15,16] = 1.000849628844304
$data[15,17] = 1.000318609048226
$data[15,18] = 1.000318609048226
$data[15,19] = 1.000113281403317
$data[16,0] = 15
$data[16,1] = "Holden5"
$data[16,2] = 1.004168461988193
$data[16,3] = 0.9989578813157612
$data[16,4] = 1.004168461988193
$data[16,5] = 0.9989578813157612
$data[16,6] = 0.9989578813157612
$data[16,7] = 1.002292648787026
$data[16,8] = 0.997221025856653
$data[16,9] = 0.9989578813157612
$data[16,10] = 0.9989578813157612
$data[16,11] = 0.9989578813157612
$data[16,12] = 1.001563171651977
$data[16,13] = 1.001563171651977
$data[16,14] = 1.001806330696993
$data[16,15] = 1.000694741539905
$data[16,16] = 1.000694741539905
$data[16,17] = 1.000260526483869
$data[16,18] = 1.000260526483869
$data[16,19] = 1.000092630096526
$data[17,0] = 16
$data[17,1] = "Holden10"
$data[17,2] = 1.002294025908786
$data[17,3] = 0.9994264917522137
$data[17,4] = 1.002294025908786
$data[17,5] = 0.9994264917522137
$data[17,6] = 0.9994264917522137
$data[17,7] = 1.001261710826346
$data[17,8] = 0.9984706496319896
$data[17,9] = 0.9994264917522137
$data[17,10] = 0.9994264917522137
$data[17,11] = 0.9994264917522137
$data[17,12] = 1.0008602588305
$data[17,13] = 1.0008602588305
$data[17,14] = 1.000994076162449
$data[17,15] = 1.000382336471071
$data[17,16] = 1.000382336471071
$data[17,17] = 1.000143375291357
$data[17,18] = 1.000143375291357
$data[17,19] = 1.000050976937294
$data[18,0] = 17
$data[18,1] = "Holden15"
$data[18,2] = 1.002501008171121
$data[18,3] = 0.9993747421608823
$data[18,4] = 1.002501008171121
$data[18,5] = 0.9993747421608823
$data[18,6] = 0.9993747421608823
$data[18,7] = 1.001375555419878
$data[18,8] = 0.9983326561189785
$data[18,9] = 0.9993747421608823
$data[18,10] = 0.9993747421608823
$data[18,11] = 0.9993747421608823
$data[18,12] = 1.000937875166002
$data[18,13] = 1.000937875166002
$data[18,14] = 1.001083768583961
$data[18,15] = 1.000416830830962
$data[18,16] = 1.000416830830962
$data[18,17] = 1.000156308663442
$data[18,18] = 1.000156308663442
$data[18,19] = 1.000055574365438
$data[19,0] = 18
$data[19,1] = "HexGrid-90degTilt2.5degRes"
$data[19,2] = 1.000001079611685
$data[19,3] = 0.999999727622689
$data[19,4] = 1.000001079611685
$data[19,5] = 0.999999727622689
$data[19,6] = 0.999999727622689
$data[19,7] = 1.000000592539203
$data[19,8] = 0.9999992795846588
$data[19,9] = 0.999999727622689
$data[19,10] = 0.999999727622689
$data[19,11] = 0.999999727622689
$data[19,12] = 1.000000403617187
$data[19,13] = 1.000000403617187
$data[19,14] = 1.000000466591192
$data[19,15] = 1.000000178285688
$data[19,16] = 1.000000178285688
$data[19,17] = 1.000000065619938
$data[19,18] = 1.000000065619938
$data[19,19] = 1.000000022433936
$data[20,0] = 19
$data[20,1] = "HexGrid-90degTilt5degRes"
$data[20,2] = 1.000022860390393
$data[20,3] = 0.9999942822547845
$data[20,4] = 1.000022860390393
$data[20,5] = 0.9999942822547845
$data[20,6] = 0.9999942822547845
$data[20,7] = 1.0000125715139
$data[20,8] = 0.9999847593117226
$data[20,9] = 0.9999942822547845
$data[20,10] = 0.9999942822547845
$data[20,11] = 0.9999942822547845
$data[20,12] = 1.000008571322589
$data[20,13] = 1.000008571322589
$data[20,14] = 1.000009904719692
$data[20,15] = 1.000003808299987
$data[20,16] = 1.000003808299987
$data[20,17] = 1.000001426788687
$data[20,18] = 1.000001426788687
$data[20,19] = 1.000000506330061
$data[21,0] = 20
$data[21,1] = "HexGrid-90degTilt10degRes"
$data[21,2] = 1.000080486315863
$data[21,3] = 0.9999798765436096
$data[21,4] = 1.000080486315863
$data[21,5] = 0.9999798765436096
$data[21,6] = 0.9999798765436096
$data[21,7] = 1.000044262496521
$data[21,8] = 0.9999463447508222
$data[21,9] = 0.9999798765436096
$data[21,10] = 0.9999798765436096
$data[21,11] = 0.9999798765436096
$data[21,12] = 1.000030181429737
$data[21,13] = 1.000030181429737
$data[21,14] = 1.000034875118665
$data[21,15] = 1.000013413134361
$data[21,16] = 1.000013413134361
$data[21,17] = 1.000005028986673
$data[21,18] = 1.000005028986673
$data[21,19] = 1.000001787199006
$data[22,0] = 21
$data[22,1] = "HexGrid-90degTilt15degRes"
$data[22,2] = 1.000188554038776
$data[22,3] = 0.9999528625992226
$data[22,4] = 1.000188554038776
$data[22,5] = 0.9999528625992226
$data[22,6] = 0.9999528625992226
$data[22,7] = 1.000103699572714
$data[22,8] = 0.9998742982819733
$data[22,9] = 0.9999528625992226
$data[22,10] = 0.9999528625992226
$data[22,11] = 0.9999528625992226
$data[22,12] = 1.000070708318999
$data[22,13] = 1.000070708318999
$data[22,14] = 1.000081705403571
$data[22,15] = 1.000031426412407
$data[22,16] = 1.000031426412407
$data[22,17] = 1.000011785459111
$data[22,18] = 1.000011785459111
$data[22,19] = 1.000004189948522
$targetRange = $ws.Range("A1").Resize(23, 20)
$targetRange.Value = $data

# --- Re-apply the bold/centered/top/bordered style to the header row (B1:T1) and index column (A2:A23) ---
$headerRange = $ws.Range("B1:T1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$indexRange = $ws.Range("A2:A23")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
$indexRange.Borders.LineStyle = 1
